$wb = $excel.ActiveWorkbook

$wsResults = $wb.Worksheets.Item("Results")
$wsMadhu = $wb.Worksheets.Item("madhu")

# Data for the "Results" sheet (columns: id, player_id, full_name, team, price, ts)
$resultsData = @(
    @(111, 70,  "Rahul Chattaraj ",     "madhu", 20, "2025-09-02 23:37:51"),
    @(112, 61,  "Abhradip Maji ",       "madhu", 20, "2025-09-02 23:38:30"),
    @(113, 49,  "Debojyoti Saha ",      "madhu", 20, "2025-09-02 23:38:38"),
    @(114, 112, "SAGEN MURMU ",         "madhu", 20, "2025-09-02 23:38:51"),
    @(115, 1,   "Anupam Ghosh ",        "madhu", 20, "2025-09-02 23:39:01"),
    @(116, 82,  "Biplab Senapati ",     "madhu", 20, "2025-09-02 23:39:10"),
    @(117, 41,  "Ambuj Mahato",         "madhu", 20, "2025-09-02 23:39:18"),
    @(118, 94,  "ARKA DEY",             "madhu", 20, "2025-09-02 23:39:27"),
    @(119, 83,  "MRIGANKA MAHATO",      "madhu", 20, "2025-09-02 23:39:38"),
    @(120, 75,  "Amit Das ",            "madhu", 20, "2025-09-02 23:39:47"),
    @(121, 51,  "AMIT MUKHERJEE ",      "madhu", 20, "2025-09-02 23:40:01"),
    @(122, 32,  "Rahul Mahato",         "madhu", 20, "2025-09-02 23:40:16"),
    @(123, 104, "Suman Mahato ",        "madhu", 20, "2025-09-02 23:40:47"),
    @(124, 99,  "Asit Kumar Mahato ",   "madhu", 20, "2025-09-02 23:40:59"),
    @(125, 33,  "Anjan Kumar mandal ",  "madhu", 20, "2025-09-02 23:47:29"),
    @(126, 37,  "Nadim Rijbi ",         "madhu", 20, "2025-09-02 23:47:40"),
    @(127, 111, "Abhay Mandal ",        "madhu", 20, "2025-09-02 23:47:51")
)

$r = 2
foreach ($row in $resultsData) {
    $wsResults.Cells.Item($r, 1).Value = $row[0]
    $wsResults.Cells.Item($r, 2).Value = $row[1]
    $wsResults.Cells.Item($r, 3).Value = $row[2]
    $wsResults.Cells.Item($r, 4).Value = $row[3]
    $wsResults.Cells.Item($r, 5).Value = $row[4]
    $wsResults.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# Data for the "madhu" sheet (columns: Player ID, FULL NAME, Price)
$madhuData = @(
    @(70,  "Rahul Chattaraj ",    20),
    @(61,  "Abhradip Maji ",      20),
    @(49,  "Debojyoti Saha ",     20),
    @(112, "SAGEN MURMU ",        20),
    @(1,   "Anupam Ghosh ",       20),
    @(82,  "Biplab Senapati ",    20),
    @(41,  "Ambuj Mahato",        20),
    @(94,  "ARKA DEY",            20),
    @(83,  "MRIGANKA MAHATO",     20),
    @(75,  "Amit Das ",           20),
    @(51,  "AMIT MUKHERJEE ",     20),
    @(32,  "Rahul Mahato",        20),
    @(104, "Suman Mahato ",       20),
    @(99,  "Asit Kumar Mahato ",  20),
    @(33,  "Anjan Kumar mandal ", 20),
    @(37,  "Nadim Rijbi ",        20),
    @(111, "Abhay Mandal ",       20)
)

$r = 2
foreach ($row in $madhuData) {
    $wsMadhu.Cells.Item($r, 1).Value = $row[0]
    $wsMadhu.Cells.Item($r, 2).Value = $row[1]
    $wsMadhu.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
